# Commit: "Actualiza Tarea para entrega 04/10/2022"
#
# Changes:
#  1. Rename worksheet "Hoja2" -> "factores"
#  2. Fix the header typo in that sheet: B1 "Sesonal " -> "Seasons"
#     (this also makes the old "Sesonal " shared-string entry unused,
#     so it disappears from sharedStrings.xml on save, shifting every
#     later shared-string index down by one — exactly what the diff shows)
#  3. Move the active selection on that sheet from G9 to F26

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Hoja2")
$ws.Name = "factores"
$ws.Range("B1").Value = "Seasons"

$ws.Activate()
$ws.Range("F26").Select()
